# 20.04.2017 tarihli Veritabani ornekleri
#
# Replace the sample student roster with a new set of rows and add a new
# "S.No." index column in A. The old row 8 disappears (6 data rows instead
# of 7), the per-row record-index column A keeps its existing 1..6 values,
# and the data cells B:G lose the thin border they used to have (only the
# header row keeps it). Column E stays a short date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the old 7th data row (row 8) completely - this shifts nothing
#    else and shrinks the used range from A1:G8 down to A1:G7.
# ---------------------------------------------------------------------
$ws.Rows("8:8").Delete() | Out-Null

# ---------------------------------------------------------------------
# 2. The data cells (B:G, rows 2-7) no longer carry the thin border the
#    old rows had - only the header row (row 1) keeps its border/fill.
# ---------------------------------------------------------------------
$ws.Range("B2:G7").Borders.LineStyle = -4142

# ---------------------------------------------------------------------
# 3. New leading "S.No." header in column A (plain style, same as the
#    numeric index cells already sitting below it in A2:A7).
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "S.No."

# ---------------------------------------------------------------------
# 4. Overwrite the student records with the new sample data. Column A
#    (1..6) is already correct and untouched.
# ---------------------------------------------------------------------
$data = @(
    @(14219022, "Ali Eren", "Sugeçmez",   35102, "YBS", "E"),
    @(15219020, "Merve",    "Gültekin",   34965, "YBS", "K"),
    @(16219014, "Sultan",   "Erdoğan",    35453, "YBS", "K"),
    @(15219505, "Mustafa",  "Özgün",      34750, "YBS", "E"),
    @(15219018, "Tansu",    "Gökçe",      35665, "YBS", "K"),
    @(14219007, "Mehmet",   "Eskicioğlu", 35065, "YBS", "E")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value2 = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value2 = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# Date column keeps a short-date display format.
$ws.Range("E2:E7").NumberFormat = "m/d/yy"

# ---------------------------------------------------------------------
# 5. Column D ("Soyadı") now needs its width best-fit like column E.
# ---------------------------------------------------------------------
$ws.Columns(4).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 6. View settings: zoom + selection moved to B2.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 205
$ws.Range("B2").Select() | Out-Null
